$wb = $excel.ActiveWorkbook

# Sheet "展览" (rId1 / sheet1.xml)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F5").Value = 5062
$wsExhibition.Range("F11").Value = 1044
$wsExhibition.Range("F14").Value = 3740
$wsExhibition.Range("F16").Value = 145
$wsExhibition.Range("F29").Value = 277

# Sheet "全部类型" (rId4 / sheet4.xml) - aggregated view, one row offset vs 展览
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 5062
$wsAll.Range("F12").Value = 1044
$wsAll.Range("F15").Value = 3740
$wsAll.Range("F17").Value = 145
$wsAll.Range("F30").Value = 277
